# Update PLC data 2025-10-13 13:49:34
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7181
$ws.Range("C3").Value = 161651
$ws.Range("C4").Value = 152660
$ws.Range("C7").Value = 5.56
$ws.Range("C8").Value = 64.5
